$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A3 value to the new user "locked_out_user" (a new shared string)
$ws.Range("A3").Value = "locked_out_user"

# Row 2 and Row 3 (A2:B3) lose their highlighted style -> revert to default style
$ws.Range("A2:B3").Style = "Normal"

# Move the active selection to C4 (no data there - just reflects last user selection)
$ws.Range("C4").Select() | Out-Null
